$wb = $excel.ActiveWorkbook

$wsFilter = $wb.Worksheets.Item("Reviews_Filter")
$wsZoom   = $wb.Worksheets.Item("Reviews_Zoom")

# --- Reviews_Filter (sheet1): update Country/State/City/Location values ---
$wsFilter.Range("C2").Value = "CA"
$wsFilter.Range("D2").Value = "ON"
$wsFilter.Range("E2").Value = "Toronto"

# --- Reviews_Zoom (sheet2): update the zoom date range ---
$wsZoom.Range("A3").Value = 1
$wsZoom.Range("E3").Value = "October"
$wsZoom.Range("B3").Value = "Septemper"
$wsZoom.Range("D3").Value = 30

# The above Value assignments clear the "quote prefix" cell style (s="2")
# that the whole row previously carried. Restore it by copying formats
# only from a sibling cell in the same row that kept its original style.
$wsZoom.Range("C3").Copy()
$wsZoom.Range("A3").PasteSpecial(-4122)
$wsZoom.Range("B3").PasteSpecial(-4122)
$wsZoom.Range("D3").PasteSpecial(-4122)
$wsZoom.Range("E3").PasteSpecial(-4122)

# Location value is longer, set it after the shared-string churn above so
# new shared strings line up in the same append order as the authored file.
$wsFilter.Range("F2").Value = "Neural Turing Tech - Primrose, 1131 Steeles Ave. West, M2R 3W8, +14164510870"

# --- Selections / active sheet ---
# Reviews_Zoom's selection moves to B3 (sheet stays inactive).
$wsZoom.Range("B3").Select()
# Reviews_Filter becomes the active sheet with F2 selected (was
# Reviews_AdvancedFilters/C2 before).
$wsFilter.Range("F2").Select()
